$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / relabel the header row for the item template:
#   A1: Item ID   (was "ID")
#   B1: SKU Code  (was "Item Code")
#   C1: Item Name (unchanged text, now 3rd column)
#   D1: Variant/ Model Name (new 4th column)
$ws.Range("A1").Value = "Item ID"
$ws.Range("B1").Value = "SKU Code"
$ws.Range("C1").Value = "Item Name"
$ws.Range("D1").Value = "Variant/ Model Name"

# D1 should carry the same header formatting as the other "dark" header
# cell (C1) rather than the default style picked up by a plain .Value write.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Column widths: C narrows slightly, D is a brand-new column.
$ws.Columns.Item(3).ColumnWidth = 46.833333333333336
$ws.Columns.Item(4).ColumnWidth = 47.5
